# Applies the "Atualizacao de bases das ligas, do dia: 03-04-2024 as 22:09" update
# to the "Belgium First Division B" sheet: a handful of adjacent match rows had their
# row-order swapped (id/date/home/away/odds all move together) in the source feed,
# a couple of already-scheduled fixtures got refreshed closing odds, and one new
# fixture (6809874) was inserted, shifting the trailing schedule rows down by one
# and appending a brand-new final row (225).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 97
$ws.Cells.Item(97, 2).Value = 6809783  # B97
$ws.Cells.Item(97, 6).Value = "Lommel"  # F97
$ws.Cells.Item(97, 7).Value = "Patro Eisden Maasmechelen"  # G97
$ws.Cells.Item(97, 8).Value = 1  # H97
$ws.Cells.Item(97, 10).Value = "A"  # J97
$ws.Cells.Item(97, 11).Value = 2.05  # K97
$ws.Cells.Item(97, 12).Value = 3.2  # L97
$ws.Cells.Item(97, 13).Value = 3.5  # M97
$ws.Cells.Item(97, 14).Value = 1.75  # N97
$ws.Cells.Item(97, 16).Value = 4.2  # P97
$ws.Cells.Item(97, 17).Value = -0.5  # Q97
$ws.Cells.Item(97, 18).Value = 1.8  # R97
$ws.Cells.Item(97, 19).Value = 2  # S97
$ws.Cells.Item(97, 20).Value = 2.25  # T97
$ws.Cells.Item(97, 21).Value = 1.8  # U97
$ws.Cells.Item(97, 22).Value = 2  # V97
$ws.Cells.Item(97, 23).Value = -1  # W97
$ws.Cells.Item(97, 25).Value = 3.2  # Y97
$ws.Cells.Item(97, 26).Value = -1  # Z97
$ws.Cells.Item(97, 27).Value = 1  # AA97
$ws.Cells.Item(97, 28).Value = 0.8  # AB97

# Row 98
$ws.Cells.Item(98, 2).Value = 6809780  # B98
$ws.Cells.Item(98, 6).Value = "WaaslandBeveren"  # F98
$ws.Cells.Item(98, 7).Value = "FCV Dender EH"  # G98
$ws.Cells.Item(98, 8).Value = 3  # H98
$ws.Cells.Item(98, 10).Value = "H"  # J98
$ws.Cells.Item(98, 11).Value = 1.95  # K98
$ws.Cells.Item(98, 12).Value = 3.6  # L98
$ws.Cells.Item(98, 13).Value = 3.3  # M98
$ws.Cells.Item(98, 14).Value = 2.15  # N98
$ws.Cells.Item(98, 16).Value = 3  # P98
$ws.Cells.Item(98, 17).Value = -0.25  # Q98
$ws.Cells.Item(98, 18).Value = 1.9  # R98
$ws.Cells.Item(98, 19).Value = 1.9  # S98
$ws.Cells.Item(98, 20).Value = 2.75  # T98
$ws.Cells.Item(98, 21).Value = 1.975  # U98
$ws.Cells.Item(98, 22).Value = 1.825  # V98
$ws.Cells.Item(98, 23).Value = 1.15  # W98
$ws.Cells.Item(98, 25).Value = -1  # Y98
$ws.Cells.Item(98, 26).Value = 0.8999999999999999  # Z98
$ws.Cells.Item(98, 27).Value = -1  # AA98
$ws.Cells.Item(98, 28).Value = 0.9750000000000001  # AB98

# Row 100
$ws.Cells.Item(100, 2).Value = 6809785  # B100
$ws.Cells.Item(100, 6).Value = "Francs Borains"  # F100
$ws.Cells.Item(100, 7).Value = "Deinze"  # G100
$ws.Cells.Item(100, 8).Value = 0  # H100
$ws.Cells.Item(100, 9).Value = 1  # I100
$ws.Cells.Item(100, 11).Value = 4.2  # K100
$ws.Cells.Item(100, 12).Value = 3.6  # L100
$ws.Cells.Item(100, 13).Value = 1.727  # M100
$ws.Cells.Item(100, 14).Value = 4  # N100
$ws.Cells.Item(100, 15).Value = 3.5  # O100
$ws.Cells.Item(100, 16).Value = 1.8  # P100
$ws.Cells.Item(100, 17).Value = 0.5  # Q100
$ws.Cells.Item(100, 18).Value = 2  # R100
$ws.Cells.Item(100, 19).Value = 1.8  # S100
$ws.Cells.Item(100, 20).Value = 2.5  # T100
$ws.Cells.Item(100, 25).Value = 0.8  # Y100
$ws.Cells.Item(100, 27).Value = 0.8  # AA100
$ws.Cells.Item(100, 28).Value = -1  # AB100
$ws.Cells.Item(100, 29).Value = 1  # AC100

# Row 101
$ws.Cells.Item(101, 2).Value = 6809782  # B101
$ws.Cells.Item(101, 6).Value = "ZulteWaregem"  # F101
$ws.Cells.Item(101, 7).Value = "Anderlecht II"  # G101
$ws.Cells.Item(101, 8).Value = 2  # H101
$ws.Cells.Item(101, 9).Value = 5  # I101
$ws.Cells.Item(101, 11).Value = 1.444  # K101
$ws.Cells.Item(101, 12).Value = 4.75  # L101
$ws.Cells.Item(101, 13).Value = 5.5  # M101
$ws.Cells.Item(101, 14).Value = 1.444  # N101
$ws.Cells.Item(101, 15).Value = 4.75  # O101
$ws.Cells.Item(101, 16).Value = 5.5  # P101
$ws.Cells.Item(101, 17).Value = -1.25  # Q101
$ws.Cells.Item(101, 18).Value = 1.95  # R101
$ws.Cells.Item(101, 19).Value = 1.85  # S101
$ws.Cells.Item(101, 20).Value = 3  # T101
$ws.Cells.Item(101, 25).Value = 4.5  # Y101
$ws.Cells.Item(101, 27).Value = 0.8500000000000001  # AA101
$ws.Cells.Item(101, 28).Value = 0.8  # AB101
$ws.Cells.Item(101, 29).Value = -1  # AC101

# Row 108
$ws.Cells.Item(108, 2).Value = 6809788  # B108
$ws.Cells.Item(108, 6).Value = "Anderlecht II"  # F108
$ws.Cells.Item(108, 7).Value = "WaaslandBeveren"  # G108
$ws.Cells.Item(108, 8).Value = 0  # H108
$ws.Cells.Item(108, 9).Value = 2  # I108
$ws.Cells.Item(108, 11).Value = 4  # K108
$ws.Cells.Item(108, 12).Value = 3.75  # L108
$ws.Cells.Item(108, 13).Value = 1.75  # M108
$ws.Cells.Item(108, 14).Value = 3.75  # N108
$ws.Cells.Item(108, 15).Value = 3.6  # O108
$ws.Cells.Item(108, 16).Value = 1.85  # P108
$ws.Cells.Item(108, 17).Value = 0.75  # Q108
$ws.Cells.Item(108, 18).Value = 1.7  # R108
$ws.Cells.Item(108, 19).Value = 2.2  # S108
$ws.Cells.Item(108, 20).Value = 3  # T108
$ws.Cells.Item(108, 25).Value = 0.8500000000000001  # Y108
$ws.Cells.Item(108, 27).Value = 1.2  # AA108
$ws.Cells.Item(108, 28).Value = -1  # AB108
$ws.Cells.Item(108, 29).Value = 0.875  # AC108

# Row 109
$ws.Cells.Item(109, 2).Value = 6809791  # B109
$ws.Cells.Item(109, 6).Value = "Deinze"  # F109
$ws.Cells.Item(109, 7).Value = "Club Brugge II"  # G109
$ws.Cells.Item(109, 8).Value = 1  # H109
$ws.Cells.Item(109, 9).Value = 3  # I109
$ws.Cells.Item(109, 11).Value = 1.55  # K109
$ws.Cells.Item(109, 12).Value = 4.333  # L109
$ws.Cells.Item(109, 13).Value = 4.75  # M109
$ws.Cells.Item(109, 14).Value = 1.4  # N109
$ws.Cells.Item(109, 15).Value = 4.75  # O109
$ws.Cells.Item(109, 16).Value = 6  # P109
$ws.Cells.Item(109, 17).Value = -1.25  # Q109
$ws.Cells.Item(109, 18).Value = 1.875  # R109
$ws.Cells.Item(109, 19).Value = 1.975  # S109
$ws.Cells.Item(109, 20).Value = 3.25  # T109
$ws.Cells.Item(109, 25).Value = 5  # Y109
$ws.Cells.Item(109, 27).Value = 0.9750000000000001  # AA109
$ws.Cells.Item(109, 28).Value = 0.9750000000000001  # AB109
$ws.Cells.Item(109, 29).Value = -1  # AC109

# Row 178
$ws.Cells.Item(178, 2).Value = 6809842  # B178
$ws.Cells.Item(178, 6).Value = "Patro Eisden Maasmechelen"  # F178
$ws.Cells.Item(178, 7).Value = "Lierse Kempenzonen"  # G178
$ws.Cells.Item(178, 8).Value = 3  # H178
$ws.Cells.Item(178, 9).Value = 0  # I178
$ws.Cells.Item(178, 10).Value = "H"  # J178
$ws.Cells.Item(178, 11).Value = 1.533  # K178
$ws.Cells.Item(178, 12).Value = 4.2  # L178
$ws.Cells.Item(178, 13).Value = 5  # M178
$ws.Cells.Item(178, 14).Value = 1.7  # N178
$ws.Cells.Item(178, 15).Value = 3.8  # O178
$ws.Cells.Item(178, 16).Value = 4.2  # P178
$ws.Cells.Item(178, 17).Value = -0.75  # Q178
$ws.Cells.Item(178, 18).Value = 1.925  # R178
$ws.Cells.Item(178, 19).Value = 1.875  # S178
$ws.Cells.Item(178, 20).Value = 2.5  # T178
$ws.Cells.Item(178, 21).Value = 1.8  # U178
$ws.Cells.Item(178, 22).Value = 2  # V178
$ws.Cells.Item(178, 23).Value = 0.7  # W178
$ws.Cells.Item(178, 24).Value = -1  # X178
$ws.Cells.Item(178, 26).Value = 0.925  # Z178
$ws.Cells.Item(178, 27).Value = -1  # AA178
$ws.Cells.Item(178, 28).Value = 0.8  # AB178

# Row 179
$ws.Cells.Item(179, 2).Value = 6809846  # B179
$ws.Cells.Item(179, 6).Value = "WaaslandBeveren"  # F179
$ws.Cells.Item(179, 7).Value = "Deinze"  # G179
$ws.Cells.Item(179, 8).Value = 2  # H179
$ws.Cells.Item(179, 9).Value = 2  # I179
$ws.Cells.Item(179, 10).Value = "D"  # J179
$ws.Cells.Item(179, 11).Value = 1.909  # K179
$ws.Cells.Item(179, 12).Value = 3.5  # L179
$ws.Cells.Item(179, 13).Value = 3.5  # M179
$ws.Cells.Item(179, 14).Value = 2.2  # N179
$ws.Cells.Item(179, 15).Value = 3.4  # O179
$ws.Cells.Item(179, 16).Value = 2.9  # P179
$ws.Cells.Item(179, 17).Value = -0.25  # Q179
$ws.Cells.Item(179, 18).Value = 1.975  # R179
$ws.Cells.Item(179, 19).Value = 1.825  # S179
$ws.Cells.Item(179, 20).Value = 2.75  # T179
$ws.Cells.Item(179, 21).Value = 1.95  # U179
$ws.Cells.Item(179, 22).Value = 1.85  # V179
$ws.Cells.Item(179, 23).Value = -1  # W179
$ws.Cells.Item(179, 24).Value = 2.4  # X179
$ws.Cells.Item(179, 26).Value = -0.5  # Z179
$ws.Cells.Item(179, 27).Value = 0.4125  # AA179
$ws.Cells.Item(179, 28).Value = 0.95  # AB179

# Row 210
$ws.Cells.Item(210, 2).Value = 6809867  # B210
$ws.Cells.Item(210, 6).Value = "Club Brugge II"  # F210
$ws.Cells.Item(210, 7).Value = "FCV Dender EH"  # G210
$ws.Cells.Item(210, 8).Value = 0  # H210
$ws.Cells.Item(210, 9).Value = 1  # I210
$ws.Cells.Item(210, 10).Value = "A"  # J210
$ws.Cells.Item(210, 11).Value = 6  # K210
$ws.Cells.Item(210, 13).Value = 1.5  # M210
$ws.Cells.Item(210, 14).Value = 6.5  # N210
$ws.Cells.Item(210, 15).Value = 4.5  # O210
$ws.Cells.Item(210, 16).Value = 1.5  # P210
$ws.Cells.Item(210, 17).Value = 1.25  # Q210
$ws.Cells.Item(210, 18).Value = 1.8  # R210
$ws.Cells.Item(210, 19).Value = 2  # S210
$ws.Cells.Item(210, 20).Value = 3  # T210
$ws.Cells.Item(210, 21).Value = 1.95  # U210
$ws.Cells.Item(210, 22).Value = 1.85  # V210
$ws.Cells.Item(210, 23).Value = -1  # W210
$ws.Cells.Item(210, 25).Value = 0.5  # Y210
$ws.Cells.Item(210, 26).Value = 0.4  # Z210
$ws.Cells.Item(210, 27).Value = -0.5  # AA210
$ws.Cells.Item(210, 29).Value = 0.8500000000000001  # AC210

# Row 211
$ws.Cells.Item(211, 2).Value = 6809869  # B211
$ws.Cells.Item(211, 6).Value = "KFCO Beerschot Wilrijk"  # F211
$ws.Cells.Item(211, 7).Value = "Genk II"  # G211
$ws.Cells.Item(211, 8).Value = 1  # H211
$ws.Cells.Item(211, 9).Value = 0  # I211
$ws.Cells.Item(211, 10).Value = "H"  # J211
$ws.Cells.Item(211, 11).Value = 1.5  # K211
$ws.Cells.Item(211, 13).Value = 6  # M211
$ws.Cells.Item(211, 14).Value = 1.333  # N211
$ws.Cells.Item(211, 15).Value = 5.25  # O211
$ws.Cells.Item(211, 16).Value = 8.5  # P211
$ws.Cells.Item(211, 17).Value = -1.5  # Q211
$ws.Cells.Item(211, 18).Value = 1.85  # R211
$ws.Cells.Item(211, 19).Value = 1.95  # S211
$ws.Cells.Item(211, 20).Value = 3.5  # T211
$ws.Cells.Item(211, 21).Value = 1.975  # U211
$ws.Cells.Item(211, 22).Value = 1.825  # V211
$ws.Cells.Item(211, 23).Value = 0.333  # W211
$ws.Cells.Item(211, 25).Value = -1  # Y211
$ws.Cells.Item(211, 26).Value = -1  # Z211
$ws.Cells.Item(211, 27).Value = 0.95  # AA211
$ws.Cells.Item(211, 29).Value = 0.825  # AC211

# Row 218
$ws.Cells.Item(218, 18).Value = 1.85  # R218
$ws.Cells.Item(218, 19).Value = 2  # S218
$ws.Cells.Item(218, 21).Value = 1.875  # U218
$ws.Cells.Item(218, 22).Value = 1.975  # V218

# Row 219
$ws.Cells.Item(219, 14).Value = 2.15  # N219
$ws.Cells.Item(219, 16).Value = 3.3  # P219

# Row 220
$ws.Cells.Item(220, 2).Value = 6809874  # B220
$ws.Cells.Item(220, 5).Value = 45388.45833333334  # E220
$ws.Cells.Item(220, 6).Value = "Lierse Kempenzonen"  # F220
$ws.Cells.Item(220, 7).Value = "Club Brugge II"  # G220
$ws.Cells.Item(220, 11).Value = 1.85  # K220
$ws.Cells.Item(220, 12).Value = 3.6  # L220
$ws.Cells.Item(220, 13).Value = 4  # M220
$ws.Cells.Item(220, 14).Value = 1.95  # N220
$ws.Cells.Item(220, 15).Value = 3.6  # O220
$ws.Cells.Item(220, 16).Value = 3.6  # P220
$ws.Cells.Item(220, 17).Value = -0.5  # Q220
$ws.Cells.Item(220, 18).Value = 1.95  # R220
$ws.Cells.Item(220, 19).Value = 1.9  # S220
$ws.Cells.Item(220, 21).Value = 1.925  # U220
$ws.Cells.Item(220, 22).Value = 1.925  # V220

# Row 221
$ws.Cells.Item(221, 2).Value = 6809871  # B221
$ws.Cells.Item(221, 6).Value = "Lommel"  # F221
$ws.Cells.Item(221, 7).Value = "Seraing United"  # G221
$ws.Cells.Item(221, 11).Value = 1.5  # K221
$ws.Cells.Item(221, 12).Value = 4  # L221
$ws.Cells.Item(221, 13).Value = 6  # M221
$ws.Cells.Item(221, 14).Value = 1.6  # N221
$ws.Cells.Item(221, 15).Value = 3.8  # O221
$ws.Cells.Item(221, 16).Value = 5.5  # P221
$ws.Cells.Item(221, 17).Value = -1  # Q221
$ws.Cells.Item(221, 18).Value = 2  # R221
$ws.Cells.Item(221, 19).Value = 1.85  # S221
$ws.Cells.Item(221, 20).Value = 3  # T221
$ws.Cells.Item(221, 21).Value = 1.95  # U221
$ws.Cells.Item(221, 22).Value = 1.9  # V221

# Row 222
$ws.Cells.Item(222, 2).Value = 6809873  # B222
$ws.Cells.Item(222, 5).Value = 45388.625  # E222
$ws.Cells.Item(222, 6).Value = "FCV Dender EH"  # F222
$ws.Cells.Item(222, 7).Value = "Deinze"  # G222
$ws.Cells.Item(222, 11).Value = 1.833  # K222
$ws.Cells.Item(222, 12).Value = 3.75  # L222
$ws.Cells.Item(222, 13).Value = 3.8  # M222
$ws.Cells.Item(222, 14).Value = 1.85  # N222
$ws.Cells.Item(222, 15).Value = 3.75  # O222
$ws.Cells.Item(222, 16).Value = 3.8  # P222
$ws.Cells.Item(222, 17).Value = -0.5  # Q222
$ws.Cells.Item(222, 18).Value = 1.9  # R222
$ws.Cells.Item(222, 19).Value = 1.95  # S222
$ws.Cells.Item(222, 20).Value = 2.75  # T222
$ws.Cells.Item(222, 21).Value = 1.85  # U222
$ws.Cells.Item(222, 22).Value = 2  # V222

# Row 223
$ws.Cells.Item(223, 2).Value = 6809876  # B223
$ws.Cells.Item(223, 5).Value = 45389.35416666666  # E223
$ws.Cells.Item(223, 6).Value = "Genk II"  # F223
$ws.Cells.Item(223, 7).Value = "Standard Liege II"  # G223
$ws.Cells.Item(223, 11).Value = 1.444  # K223
$ws.Cells.Item(223, 12).Value = 4.75  # L223
$ws.Cells.Item(223, 13).Value = 6  # M223
$ws.Cells.Item(223, 14).Value = 1.444  # N223
$ws.Cells.Item(223, 15).Value = 4.75  # O223
$ws.Cells.Item(223, 16).Value = 6  # P223
$ws.Cells.Item(223, 17).Value = -1.25  # Q223
$ws.Cells.Item(223, 18).Value = 2  # R223
$ws.Cells.Item(223, 19).Value = 1.85  # S223
$ws.Cells.Item(223, 20).Value = 3.25  # T223
$ws.Cells.Item(223, 21).Value = 1.95  # U223
$ws.Cells.Item(223, 22).Value = 1.9  # V223

# Row 224
$ws.Cells.Item(224, 2).Value = 6811668  # B224
$ws.Cells.Item(224, 5).Value = 45389.45833333334  # E224
$ws.Cells.Item(224, 6).Value = "KV Oostende"  # F224
$ws.Cells.Item(224, 7).Value = "KFCO Beerschot Wilrijk"  # G224
$ws.Cells.Item(224, 11).Value = 4  # K224
$ws.Cells.Item(224, 12).Value = 3.75  # L224
$ws.Cells.Item(224, 13).Value = 1.833  # M224
$ws.Cells.Item(224, 14).Value = 4  # N224
$ws.Cells.Item(224, 15).Value = 3.75  # O224
$ws.Cells.Item(224, 16).Value = 1.833  # P224
$ws.Cells.Item(224, 17).Value = 0.5  # Q224
$ws.Cells.Item(224, 18).Value = 1.975  # R224
$ws.Cells.Item(224, 19).Value = 1.875  # S224

# Row 225
$ws.Cells.Item(225, 1).Value = 223  # A225
$ws.Cells.Item(225, 2).Value = 6811669  # B225
$ws.Cells.Item(225, 3).Value = "Belgium First Division B"  # C225
$ws.Cells.Item(225, 4).Value = "Belgium First Division B"  # D225
$ws.Cells.Item(225, 5).Value = 45389.59375  # E225
$ws.Cells.Item(225, 6).Value = "Francs Borains"  # F225
$ws.Cells.Item(225, 7).Value = "RFC Liege"  # G225
$ws.Cells.Item(225, 11).Value = 2.8  # K225
$ws.Cells.Item(225, 12).Value = 3.4  # L225
$ws.Cells.Item(225, 13).Value = 2.45  # M225
$ws.Cells.Item(225, 14).Value = 2.8  # N225
$ws.Cells.Item(225, 15).Value = 3.4  # O225
$ws.Cells.Item(225, 16).Value = 2.45  # P225
$ws.Cells.Item(225, 17).Value = 0  # Q225
$ws.Cells.Item(225, 18).Value = 2.1  # R225
$ws.Cells.Item(225, 19).Value = 1.775  # S225
$ws.Cells.Item(225, 20).Value = 2.75  # T225
$ws.Cells.Item(225, 21).Value = 1.95  # U225
$ws.Cells.Item(225, 22).Value = 1.9  # V225
$ws.Cells.Item(225, 23).Value = 0  # W225
$ws.Cells.Item(225, 24).Value = 0  # X225
$ws.Cells.Item(225, 25).Value = 0  # Y225
$ws.Cells.Item(225, 26).Value = 0  # Z225
$ws.Cells.Item(225, 27).Value = 0  # AA225

# Row 225 is brand new; give its id (A) and date (E) cells the same number formats/
# styles used by every other data row (style index 1 = bold/centered/bordered "id"
# column, style index 2 = the YYYY-MM-DD HH:MM:SS date format) by copying them from
# the row directly above.
$ws.Range("A224").Copy()
$ws.Range("A225").PasteSpecial(-4122)
$ws.Range("E224").Copy()
$ws.Range("E225").PasteSpecial(-4122)
$excel.CutCopyMode = $false
